$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.829.39"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.72"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.73"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("E9").Value = "  -3.92%  "
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.380"
$ws.Range("E11").Value = "  -3.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.58"
$ws.Range("E12").Value = "  -17.34%  "
$ws.Range("D13").Value = "3.235.05"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.36"
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("D15").Value = "63.460.24"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("E16").Value = "  -3.55%  "
$ws.Range("D17").Value = "2.751.19"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.08"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.78"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.59"
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.70"
$ws.Range("E21").Value = "  -4.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.532"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.95"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.38"
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.92"
$ws.Range("E29").Value = "  -4.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.99"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.74"
$ws.Range("E32").Value = "  -3.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.10"
$ws.Range("E33").Value = "  -2.71%  "
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.974"
$ws.Range("E38").Value = "  -4.00%  "
$ws.Range("E39").Value = "  +4.42%  "
$ws.Range("E40").Value = "  -4.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "328.19"
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("E42").Value = "  -1.21%  "
$ws.Range("E43").Value = "  -3.37%  "
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.27"
$ws.Range("E45").Value = "  -4.00%  "
$ws.Range("E46").Value = "  -2.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "134.66"
$ws.Range("E47").Value = "  -2.69%  "
$ws.Range("E48").Value = "  -4.71%  "
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  +0.26%  "
